# Automatische test-sync: 2025-06-22 22:11:50
# Adds the new "Status van mijn bestelling" log entry as row 57 on the
# "Logs" sheet, extends the two conditional-formatting ranges that were
# sized for rows 2:56 to cover the new row 2:57, and re-orders three
# category rows (8-10) on the "Dashboard" sheet to match the refreshed
# pivot-style summary (Bestelling/Levering moves up, Factuur/Administratie
# moves down, and its count increments from 3 to 4).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row 57
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(57, 1).Value = "Status van mijn bestelling"
$logs.Cells.Item(57, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(57, 3).Value = "Wanneer wordt mijn bestelling bezorgd?"
$logs.Cells.Item(57, 4).Value = "Bestelling / Levering"
$logs.Cells.Item(57, 5).Value = "Beste klant,`nBedankt voor je e-mail. Om de status van je bestelling te controleren en meer te weten te komen over de verwachte bezorgdatum, heb ik je bestelnummer nodig. Zou je zo vriendelijk willen zijn om je bestelnummer met mij te delen, zodat ik dit verder voor je kan nakijken?`nIk hoor graag van je om je verder te helpen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item(57, 6).Value = "2025-06-22 22:11:31"
$logs.Cells.Item(57, 7).Value = "Ja"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend the conditional-formatting sqref ranges so they
#    cover the freshly added row (D2:D56 -> D2:D57, G2:G56 -> G2:G57)
# ---------------------------------------------------------------------
$catFormats = $logs.Range("D2:D56").FormatConditions
$catFormats.Item(1).ModifyAppliesToRange($logs.Range("D2:D57"))

$answeredFormats = $logs.Range("G2:G56").FormatConditions
$answeredFormats.Item(1).ModifyAppliesToRange($logs.Range("G2:G57"))

# ---------------------------------------------------------------------
# 3. Dashboard sheet: re-order the category summary rows 8-10
#    Before: 8=Samenwerking/Partnerverzoek(4) 9=Factuur/Administratie(3) 10=Bestelling/Levering(3)
#    After : 8=Bestelling/Levering(4)         9=Samenwerking/Partnerverzoek(4) 10=Factuur/Administratie(3)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(8, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(8, 2).Value = 4

$dash.Cells.Item(9, 1).Value = "Samenwerking / Partnerverzoek"
$dash.Cells.Item(9, 2).Value = 4

$dash.Cells.Item(10, 1).Value = "Factuur / Administratie"
$dash.Cells.Item(10, 2).Value = 3
